# "updated on Feb 08 night" - append 7 new COVID-19 case rows (35-41) to the
# line-list, plus two trailing blank-but-styled rows (42-43), and move the
# selection to J42 (matching a Feb-08-night data refresh of rawdata.xlsx).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New records, in the same column order as the header row:
# Case, lat, lon, Date, Age, Gender, Home, Stay, Visited, Related, Status
$newRows = @(
    @{ Case=34; Lat=1.355904;          Lon=103.838301;          Date="Feb-07"; Age=40; Gender="Female"; Home="Singapore"; Stay="Sin Ming Road";          Visited="Yong Thai Hang, GP Clinic, NCID"; Related="19" },
    @{ Case=35; Lat=1.288771;          Lon=103.821619;          Date="Feb-07"; Age=64; Gender="Male";   Home="Singapore"; Stay="Henderson Crescent";     Visited="Bukit Merah Polyclinic, SGH, Redhill Market, Hawker centre at Bukit Merah"; Related="" },
    @{ Case=36; Lat=1.359237;          Lon=103.751117;          Date="Feb-07"; Age=38; Gender="Female"; Home="Singapore"; Stay="Bukit Batok Street 31";  Visited="Grand Hyatt Singapore, Johor Bahru, GP clinic, NCID, KK Women’s and Children’s Hospital, Ng Teng Fong General Hospital"; Related="30" },
    @{ Case=37; Lat=1.347247;          Lon=103.733663;          Date="Feb-07"; Age=53; Gender="Male";   Home="Singapore"; Stay="Jurong East Street 32";  Visited="GP clinics, Ng Teng Fong General Hospital, NCID"; Related="" },
    @{ Case=38; Lat=1.380107;          Lon=103.741224;          Date="Feb-08"; Age=52; Gender="Female"; Home="Singapore"; Stay="Choa Chu Kang Avenue 3"; Visited="Choa Chu Kang Polyclinic, The Life Church and Missions Singapore, Marina Bay Sands, Chinatown and Plaza Singapura"; Related="" },
    @{ Case=39; Lat=1.342497;          Lon=103.705433;          Date="Feb-08"; Age=51; Gender="Male";   Home="Singapore"; Stay="Jurong West Central";    Visited="Malaysia, GP clinics, NCID, Grand Hyatt Singapore"; Related="30" },
    @{ Case=40; Lat=1.333001;          Lon=103.926856;          Date="Feb-08"; Age=36; Gender="Male";   Home="Singapore"; Stay="Bedok North Street";     Visited="GP clinic, Yong Thai Hang"; Related="19" }
)

$startRow = 35

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rec = $newRows[$i]

    # Pull formatting (styles, text-as-date format on D, right-aligned J, etc.)
    # from the row directly above so every new row matches the existing ones.
    $ws.Range("A" + ($r - 1) + ":K" + ($r - 1)).Copy()
    $ws.Range("A" + $r + ":K" + $r).PasteSpecial(-4122)
    $ws.Rows.Item($r).RowHeight = 21

    $ws.Cells.Item($r, 1).Value = $rec.Case
    $ws.Cells.Item($r, 2).Value = $rec.Lat
    $ws.Cells.Item($r, 3).Value = $rec.Lon
    $ws.Cells.Item($r, 4).Value = $rec.Date
    $ws.Cells.Item($r, 5).Value = $rec.Age
    $ws.Cells.Item($r, 6).Value = $rec.Gender
    $ws.Cells.Item($r, 7).Value = $rec.Home
    $ws.Cells.Item($r, 8).Value = $rec.Stay
    $ws.Cells.Item($r, 9).Value = $rec.Visited
    if ($rec.Related -ne "") {
        $ws.Cells.Item($r, 10).Value = $rec.Related
    }
}

# Two trailing rows that only carry column A's style (no values at all), same
# as how the sheet ended before.
foreach ($r in 42, 43) {
    $ws.Range("A" + ($r - 1)).Copy()
    $ws.Range("A" + $r).PasteSpecial(-4122)
    $ws.Rows.Item($r).RowHeight = 21
}

# Scroll the view down and land the selection on J42, like the author did
# after typing in the last related-case value.
$excel.ActiveWindow.ScrollRow = 21
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("J42").Select()
